$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A1").Value = "Q_no"
$ws.Range("B1").Value = "Question"
$ws.Range("C1").Value = "Option1"
$ws.Range("D1").Value = "Option2"
$ws.Range("E1").Value = "Option3"
$ws.Range("F1").Value = "Option4"
$ws.Range("G1").Value = "Solution"
$ws.Range("H1").Value = "Remarks"

for ($i = 1; $i -le 10; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $i
    $ws.Cells.Item($row, 2).Value = "Question $i"
    $ws.Cells.Item($row, 3).Value = "Option1"
    $ws.Cells.Item($row, 4).Value = "Option2"
    $ws.Cells.Item($row, 5).Value = "Option3"
    $ws.Cells.Item($row, 6).Value = "Option4"
    $ws.Cells.Item($row, 7).Value = 3
    $ws.Cells.Item($row, 8).Value = "Op3 is right because.."
}

$ws.Range("H3").Select()
